# localization-status.xlsx -- "Generate Report for Archive"
#
# 1. All cells whose status is "Ready for handoff" move to "In Translation".
#    (There is a single shared string used by every status cell across the
#    three sheets, so every occurrence must be rewritten for the engine to
#    collapse back down to one shared-string entry, matching the source
#    diff which edits the <si> text in place.)
# 2. The "zh-cn"/"de-de" status columns on the Overview sheet (E & F) and
#    the "Status" column (C) on the per-locale sheets shrink from the old
#    width to the new, narrower width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1. Status text: "Ready for handoff" -> "In Translation" -------------

$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- 2. Narrow the status columns -----------------------------------------
# Old stored width ~17.216 chars, new stored width ~13.410 chars; that is a
# ColumnWidth (COM, "chars minus the 5/6 gridline padding") of ~12.5.
# NOTE: use the numeric column index with Columns.Item(...) -- letter
# indices (e.g. Columns.Item("E")) are not resolved correctly here.

$newColumnWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newColumnWidth   # column E
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth   # column F

$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth        # column C

$dede.Columns.Item(3).ColumnWidth = $newColumnWidth        # column C
